# Update odds values in Sheet1 for the 2025-05-29 FlashScore weekly games
# workbook, per the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5 (Antwerp - Charleroi) ---
$ws.Range("G5").Value = 2.38
$ws.Range("I5").Value = 2.7

# --- Row 7 (Jaguares de Cordoba - Inter Palmira) ---
$ws.Range("G7").Value = 1.65
$ws.Range("H7").Value = 3.55
$ws.Range("I7").Value = 4.85
$ws.Range("R7").Value = 1.88
$ws.Range("S7").Value = 1.72
$ws.Range("T7").Value = 6.2
$ws.Range("U7").Value = 7.4
$ws.Range("W7").Value = 12.5
$ws.Range("Y7").Value = 30
$ws.Range("AA7").Value = 7
$ws.Range("AB7").Value = 17.5
$ws.Range("AE7").Value = 11.75
$ws.Range("AF7").Value = 27
$ws.Range("AG7").Value = 16
$ws.Range("AH7").Value = 90
$ws.Range("AI7").Value = 55
$ws.Range("AJ7").Value = 60

# --- Row 23 (Vestmannaeyjar - Hafnarfjordur) ---
$ws.Range("L23").Value = 1.16
$ws.Range("M23").Value = 4.55
$ws.Range("N23").Value = 1.5
$ws.Range("O23").Value = 2.4
$ws.Range("P23").Value = 1.28
$ws.Range("Q23").Value = 3.35
$ws.Range("R23").Value = 1.44
$ws.Range("S23").Value = 2.6
$ws.Range("T23").Value = 13.5
$ws.Range("U23").Value = 17
$ws.Range("W23").Value = 32
$ws.Range("AA23").Value = 7.5
$ws.Range("AB23").Value = 10.75
$ws.Range("AC23").Value = 32
$ws.Range("AD23").Value = 175
$ws.Range("AE23").Value = 12.5
$ws.Range("AF23").Value = 15.5
$ws.Range("AG23").Value = 9.5
$ws.Range("AH23").Value = 28
$ws.Range("AI23").Value = 17.5
$ws.Range("AJ23").Value = 20

# --- Row 25 (Cremonese - Spezia) ---
$ws.Range("G25").Value = 2.5
$ws.Range("I25").Value = 3
$ws.Range("J25").Value = 1.08
$ws.Range("K25").Value = 8
$ws.Range("L25").Value = 1.4
$ws.Range("M25").Value = 2.75
$ws.Range("N25").Value = 2.25
$ws.Range("O25").Value = 1.62
$ws.Range("R25").Value = 1.91
$ws.Range("S25").Value = 1.91
$ws.Range("U25").Value = 12
$ws.Range("V25").Value = 10
$ws.Range("Z25").Value = 8
$ws.Range("AE25").Value = 8.5
$ws.Range("AF25").Value = 13
$ws.Range("AG25").Value = 11
$ws.Range("AH25").Value = 29

# --- Row 34 (Karagumruk - Bandirmaspor) ---
$ws.Range("N34").Value = 1.88
$ws.Range("O34").Value = 1.93
